$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (91) to the data dictionary describing the "target" column.
$ws.Range("A91").Value = "x02"
$ws.Range("B91").Value = "target"
$ws.Range("C91").Value = $false
$ws.Range("D91").Value = $false
$ws.Range("E91").Value = $false
$ws.Range("F91").Value = "numeric"
$ws.Range("G91").Value = "categorical"
$ws.Range("H91").Value = "float64"
$ws.Range("J91").Value = "n/a"
$ws.Range("K91").Value = "n/a"
$ws.Range("L91").Value = "User has history of conferences and/or favors attending conferences to learn"

# Match the author's final view state: new last row selected.
$ws.Activate()
$ws.Range("A91").Select()
